$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (avoid numeric auto-conversion) for the Price/Volume columns
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "57.882.04"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "2.447.55"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").Value = "510.97"
$ws.Range("E5").Value = "  -1.80%  "

# Row 6
$ws.Range("D6").Value = "129.89"
$ws.Range("E6").Value = "  -1.17%  "

# Row 7
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.15%  "

# Row 8
$ws.Range("D8").Value = "0.550"

# Row 9
$ws.Range("D9").Value = "2.467.70"
$ws.Range("E9").Value = "  -0.42%  "

# Row 10
$ws.Range("D10").Value = "0.0961"
$ws.Range("E10").Value = "  -3.31%  "

# Row 11
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  +0.03%  "

# Row 12
$ws.Range("D12").Value = "5.20"
$ws.Range("E12").Value = "  -2.89%  "

# Row 13
$ws.Range("D13").Value = "0.328"
$ws.Range("E13").Value = "  -4.56%  "

# Row 14
$ws.Range("D14").Value = "2.887.55"
$ws.Range("E14").Value = "  -0.91%  "

# Row 15
$ws.Range("D15").Value = "57.814.92"
$ws.Range("E15").Value = "  -0.44%  "

# Row 16
$ws.Range("D16").Value = "21.88"
$ws.Range("E16").Value = "  -0.89%  "

# Row 17
$ws.Range("D17").Value = "0.0000133"
$ws.Range("E17").Value = "  -2.23%  "

# Row 18
$ws.Range("D18").Value = "2.463.76"
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
$ws.Range("D19").Value = "10.52"
$ws.Range("E19").Value = "  -2.98%  "

# Row 20
$ws.Range("D20").Value = "318.22"
$ws.Range("E20").Value = "  -0.32%  "

# Row 21
$ws.Range("D21").Value = "4.12"
$ws.Range("E21").Value = "  -1.27%  "

# Row 22
$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.13%  "

# Row 23
$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").Value = "5.93"
$ws.Range("E23").Value = "  +3.23%  "

# Row 24
$ws.Range("D24").Value = "63.27"
$ws.Range("E24").Value = "  -1.32%  "

# Row 25
$ws.Range("D25").Value = "0.401"
$ws.Range("E25").Value = "  -2.10%  "

# Row 26
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.57%  "

# Row 27
$ws.Range("E27").Value = "  -0.31%  "

# Row 28
$ws.Range("D28").Value = "7.25"
$ws.Range("E28").Value = "  -1.65%  "

# Row 29
$ws.Range("D29").Value = "168.13"
$ws.Range("E29").Value = "  +0.76%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0729"
$ws.Range("E30").Value = "  -3.50%  "

# Row 31
$ws.Range("D31").Value = "1.66"
$ws.Range("E31").Value = "  -2.24%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  -1.28%  "

# Row 33
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "6.15"
$ws.Range("E33").Value = "  -2.75%  "

# Row 34
$ws.Range("D34").Value = "0.998"
$ws.Range("E34").Value = "  -0.07%  "

# Row 35
$ws.Range("E35").Value = "  -0.21%  "

# Row 36
$ws.Range("D36").Value = "17.76"
$ws.Range("E36").Value = "  -1.93%  "

# Row 37
$ws.Range("D37").Value = "1.27"
$ws.Range("E37").Value = "  -3.77%  "

# Row 38
$ws.Range("D38").Value = "3.90"
$ws.Range("E38").Value = "  -2.15%  "

# Row 39
$ws.Range("D39").Value = "36.55"
$ws.Range("E39").Value = "  -0.14%  "

# Row 40
$ws.Range("D40").Value = "1.45"
$ws.Range("E40").Value = "  -2.19%  "

# Row 41
$ws.Range("D41").Value = "0.759"
$ws.Range("E41").Value = "  -4.26%  "

# Row 42
$ws.Range("D42").Value = "270.35"
$ws.Range("E42").Value = "  -2.05%  "

# Row 43
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D43").Value = "5.00"
$ws.Range("E43").Value = "  -0.85%  "

# Row 44
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "3.38"
$ws.Range("E44").Value = "  -2.98%  "

# Row 45
$ws.Range("D45").Value = "0.586"
$ws.Range("E45").Value = "  -1.41%  "

# Row 46
$ws.Range("D46").Value = "0.0913"
$ws.Range("E46").Value = "  +0.78%  "

# Row 47
$ws.Range("D47").Value = "120.35"
$ws.Range("E47").Value = "  -4.69%  "

# Row 48
$ws.Range("D48").Value = "0.0488"
$ws.Range("E48").Value = "  -0.42%  "

# Row 49
$ws.Range("D49").Value = "17.28"
$ws.Range("E49").Value = "  -3.60%  "

# Row 50
$ws.Range("E50").Value = "  -2.25%  "

# Row 51
$ws.Range("D51").Value = "16.69"
$ws.Range("E51").Value = "  -2.68%  "

# Remove the temporary text-number-format so styling matches original (no explicit style)
$ws.Range("D2:E51").ClearFormats()
